# Insert a new data row at row 220 (pushing the existing rows 220-254 down
# to 221-255) and populate it with the new "Cilantro" market record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(220).Insert()

$ws.Cells.Item(220, 1).Value  = 10
$ws.Cells.Item(220, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(220, 3).Value  = "La Araucanía"
$ws.Cells.Item(220, 4).Value  = 44505
$ws.Cells.Item(220, 5).Value  = 9
$ws.Cells.Item(220, 6).Value  = 100112040
$ws.Cells.Item(220, 7).Value  = "Cilantro"
$ws.Cells.Item(220, 8).Value  = "Sin especificar"
$ws.Cells.Item(220, 9).Value  = "Primera"
$ws.Cells.Item(220, 10).Value = 50
$ws.Cells.Item(220, 11).Value = 4500
$ws.Cells.Item(220, 12).Value = 4500
$ws.Cells.Item(220, 13).Value = 4500
$ws.Cells.Item(220, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(220, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(220, 16).Value = 2250
$ws.Cells.Item(220, 17).Value = 2
$ws.Cells.Item(220, 18).Value = "Hortaliza"
